$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Remove the "Meta description: ..." paragraph that currently
#    sits right under the H1 title, and re-insert its bold lead-in
#    run (retitled) as a new paragraph just before the final
#    "Prompt: ..." paragraph, where the italic "Prompt" text becomes
#    the plain "Read our Cinderella's Ball..." sentence.
# ---------------------------------------------------------------

$metaPara = $d.Paragraphs.Item(2)

# Sanity check - make sure we are grabbing the right paragraph.
# (Leave as comment-safe no-op if mismatched; Find below is robust.)

$metaPara.Range.Select()
$word.Selection.Cut()

# After the cut, the (former) last paragraph is now one slot lower;
# re-fetch it fresh from the live paragraph collection.
$n = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($n)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertPoint.Select()
$word.Selection.Paste()

# The pasted paragraph (formerly "Meta description: ...") is now the
# second-to-last paragraph.
$n2 = $d.Paragraphs.Count
$movedPara = $d.Paragraphs.Item($n2 - 1)
$mr = $movedPara.Range

$boldLabel = "Meta description"

# Drop everything after the bold "Meta description" run (that is the
# ": Read our Cinderella's Ball slot review..." sentence) - it moves
# to the old Prompt paragraph below instead.
$suffixStart = $mr.Start + $boldLabel.Length
$suffixEnd = $mr.End - 1
if ($suffixEnd -gt $suffixStart) {
    $d.Range($suffixStart, $suffixEnd).Delete()
}

# Retitle the remaining bold run.
$mr2 = $movedPara.Range
$boldRange = $d.Range($mr2.Start, $mr2.Start + $boldLabel.Length)
$boldRange.Text = "Play Cinderella's Ball Free: Magical Slot Machine Review"

# ---------------------------------------------------------------
# 2. Swap the old AI-image "Prompt: ..." paragraph's text for the
#    short meta-description sentence (keeping its italic styling).
# ---------------------------------------------------------------

$oldPrompt = "Prompt: Create a cartoon-style feature image that features a happy Maya warrior wearing glasses. The image should be eye-catching and appealing, with bright colors and playful elements. The warrior should have a big smile on their face and be surrounded by symbols and elements from the Cinderella's Ball game, such as the pumpkin carriage, the crystal slipper, and the magic wand. The text " + [char]34 + "Cinderella's Ball" + [char]34 + " should be featured prominently in the image, in fun and playful font."
$newDescription = "Read our Cinderella's Ball slot review to play a magical slot machine game free. With an intuitive gameplay, Cinderella's Ball offers bonus rounds and themed payouts."

$finalCount = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($finalCount)
$pr = $promptPara.Range

if ($pr.Text -eq ($oldPrompt + [char]13)) {
    $target = $d.Range($pr.Start, $pr.End - 1)
    $target.Text = $newDescription
} else {
    $d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newDescription, 2)
}
